# Updates cryptocurrency price (D) and 1h volume change (E) columns
# Each value is prefixed with a literal apostrophe so Excel stores it
# as text (matching the original inlineStr cell type) instead of
# auto-converting number-looking strings (e.g. "0.9989") into numerics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.299.87"
$ws.Range("E2").Value = "'  -0.68%  "
$ws.Range("D3").Value = "'1.620.18"
$ws.Range("E3").Value = "'  -0.67%  "
$ws.Range("E4").Value = "'  -0.07%  "
$ws.Range("D5").Value = "'0.9989"
$ws.Range("E5").Value = "'  -0.15%  "
$ws.Range("D6").Value = "'302.26"
$ws.Range("E6").Value = "'  -0.88%  "
$ws.Range("D7").Value = "'0.3742"
$ws.Range("E7").Value = "'  +0.30%  "
$ws.Range("D8").Value = "'51.82"
$ws.Range("E8").Value = "'  +0.31%  "
$ws.Range("D9").Value = "'0.3551"
$ws.Range("E9").Value = "'  -2.72%  "
$ws.Range("D10").Value = "'0.08148"
$ws.Range("E10").Value = "'  -0.41%  "
$ws.Range("D11").Value = "'1.222"
$ws.Range("E11").Value = "'  -0.36%  "
$ws.Range("D12").Value = "'0.9992"
$ws.Range("E12").Value = "'  -0.04%  "
$ws.Range("D13").Value = "'22.14"
$ws.Range("E13").Value = "'  -1.80%  "
$ws.Range("D14").Value = "'6.450"
$ws.Range("E14").Value = "'  -1.55%  "
$ws.Range("D15").Value = "'7.259"
$ws.Range("E15").Value = "'  -0.16%  "
$ws.Range("D16").Value = "'0.00001219"
$ws.Range("E16").Value = "'  -2.53%  "
$ws.Range("D17").Value = "'1.612.63"
$ws.Range("E17").Value = "'  -1.13%  "
$ws.Range("D18").Value = "'95.23"
$ws.Range("E18").Value = "'  +0.83%  "
$ws.Range("D19").Value = "'0.06926"
$ws.Range("E19").Value = "'  -0.66%  "
$ws.Range("D20").Value = "'6.633"
$ws.Range("E20").Value = "'  +2.68%  "
$ws.Range("D21").Value = "'17.36"
$ws.Range("E21").Value = "'  -2.35%  "
$ws.Range("D22").Value = "'0.9990"
$ws.Range("E22").Value = "'  -0.17%  "
$ws.Range("D23").Value = "'12.40"
$ws.Range("E23").Value = "'  -2.85%  "
$ws.Range("D24").Value = "'23.295.77"
$ws.Range("E24").Value = "'  -0.67%  "
$ws.Range("D25").Value = "'2.515"
$ws.Range("E25").Value = "'  +2.00%  "
$ws.Range("D26").Value = "'3.084"
$ws.Range("E26").Value = "'  -3.54%  "
$ws.Range("D27").Value = "'20.95"
$ws.Range("E27").Value = "'  -2.48%  "
$ws.Range("D28").Value = "'152.71"
$ws.Range("E28").Value = "'  +1.49%  "
$ws.Range("D29").Value = "'5.172"
$ws.Range("E29").Value = "'  -3.28%  "
$ws.Range("D30").Value = "'133.24"
$ws.Range("E30").Value = "'  -0.99%  "
$ws.Range("D31").Value = "'1.793.95"
$ws.Range("E31").Value = "'  -1.06%  "
$ws.Range("D32").Value = "'1.094"
$ws.Range("E32").Value = "'  +6.93%  "
$ws.Range("D33").Value = "'6.533"
$ws.Range("E33").Value = "'  -4.44%  "
$ws.Range("D34").Value = "'11.69"
$ws.Range("E34").Value = "'  +6.39%  "
$ws.Range("D35").Value = "'2.027"
$ws.Range("E35").Value = "'  -10.56%  "
$ws.Range("D36").Value = "'0.02730"
$ws.Range("E36").Value = "'  -1.92%  "
$ws.Range("D37").Value = "'0.08724"
$ws.Range("E37").Value = "'  -0.48%  "
$ws.Range("D38").Value = "'0.2463"
$ws.Range("E38").Value = "'  -2.73%  "
$ws.Range("D39").Value = "'0.06928"
$ws.Range("E39").Value = "'  -2.87%  "
$ws.Range("D40").Value = "'5.882"
$ws.Range("E40").Value = "'  -3.11%  "
$ws.Range("D41").Value = "'12.53"
$ws.Range("E41").Value = "'  +1.73%  "
$ws.Range("D42").Value = "'0.6907"
$ws.Range("E42").Value = "'  -2.10%  "
$ws.Range("D43").Value = "'1.325"
$ws.Range("E43").Value = "'  -1.79%  "
$ws.Range("E44").Value = "'  -4.58%  "
$ws.Range("D45").Value = "'0.9989"
$ws.Range("E45").Value = "'  -0.08%  "
$ws.Range("D46").Value = "'0.6363"
$ws.Range("E46").Value = "'  -2.52%  "
$ws.Range("D47").Value = "'2.260"
$ws.Range("E47").Value = "'  -3.21%  "
$ws.Range("D48").Value = "'3.940"
$ws.Range("E48").Value = "'  -1.35%  "
$ws.Range("D49").Value = "'0.07881"
$ws.Range("E49").Value = "'  -1.89%  "
$ws.Range("D50").Value = "'126.94"
$ws.Range("E50").Value = "'  +1.44%  "
$ws.Range("D51").Value = "'1.170"
$ws.Range("E51").Value = "'  -3.01%  "
